# HC Services.xlsx update: add Haiipfy raise info to the Private sheet
# and make Private the active/selected sheet (matches author's edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Private")

# New entries below the existing WebMD/KKR/2.8B row (row 4):
# Row 2 gets a "Raise"/"Round" mini-header in columns E:F,
# Row 5 gets the Haiipfy company row (B) plus its raise amount/round (E:F).
# Written in this order so the new shared-string table entries line up
# with the source workbook (Haiipfy, Raise, 73m, Round, D).
$ws.Range("B5").Value = "Haiipfy"
$ws.Range("E2").Value = "Raise"
$ws.Range("E5").Value = "73m"
$ws.Range("F2").Value = "Round"
$ws.Range("F5").Value = "D"

# Make "Private" the active sheet/tab, with the active cell parked at F6
# (just below the newly added data), matching the saved view state.
$ws.Activate()
$ws.Range("F6").Select()
